# "Generate Report for Archive"
#
# 1. Every cell whose displayed text is "Ready for handoff" becomes
#    "In Translation" (the shared string used by the Overview sheet's
#    zh-cn/de-de status columns and by the per-language "Status" column
#    on the zh-cn / de-de sheets).
# 2. The "Status"-column width is narrowed (it had been sized to fit the
#    longer "Ready for handoff" text) on all three sheets:
#      - Overview!E:F  (zh-cn / de-de status columns)
#      - zh-cn!C       (Status column)
#      - de-de!C       (Status column)

$wb = $excel.ActiveWorkbook

# --- 1. Replace the status text everywhere it appears ------------------
# NOTE: cast Range.Text to [string] explicitly before comparing - some
# cells (e.g. the boolean "True"/"False" ones) come back as a typed
# Boolean, and comparing that un-cast against a string literal coerces
# the *string* to a bool instead, causing false positives.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $text = [string]$cell.Text
            if ($text -eq "Ready for handoff") {
                $cell.Value = "In Translation"
            }
        }
    }
}

# --- 2. Shrink the Status columns ---------------------------------------
$targetWidth = 13.4101845877511

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1:F1").ColumnWidth = $targetWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").ColumnWidth = $targetWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").ColumnWidth = $targetWidth
